# Bump schema version to 0.2.0 and new artefacts
# Target sheet: "IntermicrobialInteraction" (not the active sheet, so look it up by name)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IntermicrobialInteraction")

# Rewrite the single header row with the three new columns inserted:
#   sequence_id           (between tax_id and evidence_type)
#   method_type            (between evidence_type and reference)
#   participant_outcomes   (between reference and id)
$headers = @("participants", "tax_id", "sequence_id", "evidence_type", "method_type", "reference", "participant_outcomes", "id", "name", "description")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# New list-based data validations for the two new enum-style columns.
$methodTypeRange = $ws.Range("E2:E1048576")
$methodTypeRange.Validation.Add(3, 1, 1, '"simulation,microscopy,cultivation,sample"')
$methodTypeRange.Validation.IgnoreBlank = $true
$methodTypeRange.Validation.InCellDropdown = $true
$methodTypeRange.Validation.ShowInput = $false
$methodTypeRange.Validation.ShowError = $false

$participantOutcomesRange = $ws.Range("G2:G1048576")
$participantOutcomesRange.Validation.Add(3, 1, 1, '"not_affected,positively,negatively,not_applicable"')
$participantOutcomesRange.Validation.IgnoreBlank = $true
$participantOutcomesRange.Validation.InCellDropdown = $true
$participantOutcomesRange.Validation.ShowInput = $false
$participantOutcomesRange.Validation.ShowError = $false
